# Generate Report for Handoff
# Adds a new localization-status row (b512c95f-9fa1-403b-a53b-5c03f44ede5c.md)
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$HYPERLINK_COLOR = 15570276   # OLE/BGR encoding of RGB(0x64,0x95,0xED) == #6495ED
$DATE_FMT = "yyyy-mm-dd HH:mm:ss"

$fileName      = "b512c95f-9fa1-403b-a53b-5c03f44ede5c.md"
$pathAndName   = "e2e\b512c95f-9fa1-403b-a53b-5c03f44ede5c.md"
$extension     = ".md"
$status        = "Ready for handoff"
$hoDate        = "2016-09-02 22:48:58"

$zhHandoffFile = "b512c95f-9fa1-403b-a53b-5c03f44ede5c.b7cc76f1e5665da7fd6316deca7b8037124fce0d.zh-cn.xlf"
$zhHandoffDate = "2016-09-02 22:48:53"

$deHandoffFile = "b512c95f-9fa1-403b-a53b-5c03f44ede5c.b7cc76f1e5665da7fd6316deca7b8037124fce0d.de-de.xlf"
$deHandoffDate = "2016-09-02 22:48:58"

$commitSha = "dc670319b220b77d08341e419b5b11ee208b135a"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$fileName"

# ---------------------------------------------------------------------------
# Overview sheet (table3) — columns: File Name, Path And Name, Extension,
# Publish URL, zh-cn, de-de, Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$rOverview = $loOverview.ListRows.Count + 1   # header occupies row 1

$wsOverview.Cells.Item($rOverview, 1).Value = $fileName

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($rOverview, 2), $mdUrl, "", "", $pathAndName) | Out-Null
$wsOverview.Cells.Item($rOverview, 2).Font.Underline = $true
$wsOverview.Cells.Item($rOverview, 2).Font.Color = $HYPERLINK_COLOR

$wsOverview.Cells.Item($rOverview, 3).Value = $extension
$wsOverview.Cells.Item($rOverview, 4).Value = "'"
$wsOverview.Cells.Item($rOverview, 5).Value = $status
$wsOverview.Cells.Item($rOverview, 6).Value = $status
$wsOverview.Cells.Item($rOverview, 7).Value = $hoDate
$wsOverview.Cells.Item($rOverview, 7).NumberFormat = $DATE_FMT

# ---------------------------------------------------------------------------
# zh-cn sheet (table1) — columns: Source File Name, File Extension, Status,
# Source Path, Priority, Content Duplicate, Latest Handoff File,
# Latest Handoff Datetime, Latest Target File, Latest Handback File,
# Latest Handback DateTime, Reference Tokens, To be localized,
# Dependency From, Has metadata, Error Detail
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$rZh = $loZh.ListRows.Count + 1

$wsZh.Hyperlinks.Add($wsZh.Cells.Item($rZh, 1), $mdUrl, "", "", $fileName) | Out-Null
$wsZh.Cells.Item($rZh, 1).Font.Underline = $true
$wsZh.Cells.Item($rZh, 1).Font.Color = $HYPERLINK_COLOR

$wsZh.Cells.Item($rZh, 2).Value = $extension
$wsZh.Cells.Item($rZh, 3).Value = $status
$wsZh.Cells.Item($rZh, 4).Value = "e2e"
$wsZh.Cells.Item($rZh, 5).Value = "ht"
$wsZh.Cells.Item($rZh, 6).Value = "'False"

$wsZh.Cells.Item($rZh, 7).Value = $zhHandoffFile
$wsZh.Cells.Item($rZh, 8).Value = $zhHandoffDate
$wsZh.Cells.Item($rZh, 8).NumberFormat = $DATE_FMT

$wsZh.Cells.Item($rZh, 9).Value = "'"
$wsZh.Cells.Item($rZh, 10).Value = "'"

$wsZh.Cells.Item($rZh, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item($rZh, 11).NumberFormat = $DATE_FMT

$wsZh.Cells.Item($rZh, 12).Value = "'"
$wsZh.Cells.Item($rZh, 13).Value = "'True"
$wsZh.Cells.Item($rZh, 14).Value = "'"
$wsZh.Cells.Item($rZh, 15).Value = "'False"
$wsZh.Cells.Item($rZh, 16).Value = "'"

# ---------------------------------------------------------------------------
# de-de sheet (table2) — same column layout as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$rDe = $loDe.ListRows.Count + 1

$wsDe.Hyperlinks.Add($wsDe.Cells.Item($rDe, 1), $mdUrl, "", "", $fileName) | Out-Null
$wsDe.Cells.Item($rDe, 1).Font.Underline = $true
$wsDe.Cells.Item($rDe, 1).Font.Color = $HYPERLINK_COLOR

$wsDe.Cells.Item($rDe, 2).Value = $extension
$wsDe.Cells.Item($rDe, 3).Value = $status
$wsDe.Cells.Item($rDe, 4).Value = "e2e"
$wsDe.Cells.Item($rDe, 5).Value = "ht"
$wsDe.Cells.Item($rDe, 6).Value = "'False"

$wsDe.Cells.Item($rDe, 7).Value = $deHandoffFile
$wsDe.Cells.Item($rDe, 8).Value = $deHandoffDate
$wsDe.Cells.Item($rDe, 8).NumberFormat = $DATE_FMT

$wsDe.Cells.Item($rDe, 9).Value = "'"
$wsDe.Cells.Item($rDe, 10).Value = "'"

$wsDe.Cells.Item($rDe, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item($rDe, 11).NumberFormat = $DATE_FMT

$wsDe.Cells.Item($rDe, 12).Value = "'"
$wsDe.Cells.Item($rDe, 13).Value = "'True"
$wsDe.Cells.Item($rDe, 14).Value = "'"
$wsDe.Cells.Item($rDe, 15).Value = "'False"
$wsDe.Cells.Item($rDe, 16).Value = "'"
